{"js": "// Replace the date line and the 26 \"NN\u00d7NN=NNNN\" answer cells with their\n// updated values. Each old value is unique in the document, so a plain\n// (case-sensitive, non-wildcard) search-and-replace per pair is safe and\n// keeps the original run formatting (font/size) intact because we only\n// replace the text inside the found range, not the run properties.\nconst replacements = [\n  [\"2025-10-24 Friday\", \"2025-10-25 Saturday\"],\n  [\"87\u00d720=1740\", \"26\u00d725=650\"],\n  [\"54\u00d776=4104\", \"92\u00d798=9016\"],\n  [\"64\u00d769=4416\", \"79\u00d736=2844\"],\n  [\"82\u00d760=4920\", \"16\u00d776=1216\"],\n  [\"73\u00d716=1168\", \"50\u00d793=4650\"],\n  [\"85\u00d770=5950\", \"39\u00d766=2574\"],\n  [\"88\u00d792=8096\", \"69\u00d788=6072\"],\n  [\"38\u00d759=2242\", \"14\u00d725=350\"],\n  [\"29\u00d712=348\", \"56\u00d779=4424\"],\n  [\"82\u00d715=1230\", \"53\u00d713=689\"],\n  [\"52\u00d744=2288\", \"14\u00d786=1204\"],\n  [\"27\u00d729=783\", \"92\u00d775=6900\"],\n  [\"24\u00d763=1512\", \"66\u00d723=1518\"],\n  [\"86\u00d733=2838\", \"80\u00d742=3360\"],\n  [\"83\u00d788=7304\", \"86\u00d759=5074\"],\n  [\"56\u00d735=1960\", \"22\u00d776=1672\"],\n  [\"28\u00d786=2408\", \"97\u00d744=4268\"],\n  [\"87\u00d793=8091\", \"63\u00d754=3402\"],\n  [\"30\u00d793=2790\", \"64\u00d733=2112\"],\n  [\"46\u00d753=2438\", \"35\u00d744=1540\"],\n  [\"33\u00d781=2673\", \"37\u00d725=925\"],\n  [\"76\u00d799=7524\", \"17\u00d771=1207\"],\n  [\"66\u00d779=5214\", \"76\u00d737=2812\"],\n  [\"47\u00d716=752\", \"11\u00d756=616\"],\n  [\"19\u00d756=1064\", \"63\u00d799=6237\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and the 26 \"NN\u00d7NN=NNNN\" answer cells with their\n# updated values. Each old value is unique in the document, so a plain\n# case-sensitive Find/Replace (no wildcards) per pair is safe, and only the\n# run's text is touched (run formatting - font/size - is left as-is by Word's\n# Find & Replace).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-24 Friday\", \"2025-10-25 Saturday\"),\n    @(\"87\u00d720=1740\", \"26\u00d725=650\"),\n    @(\"54\u00d776=4104\", \"92\u00d798=9016\"),\n    @(\"64\u00d769=4416\", \"79\u00d736=2844\"),\n    @(\"82\u00d760=4920\", \"16\u00d776=1216\"),\n    @(\"73\u00d716=1168\", \"50\u00d793=4650\"),\n    @(\"85\u00d770=5950\", \"39\u00d766=2574\"),\n    @(\"88\u00d792=8096\", \"69\u00d788=6072\"),\n    @(\"38\u00d759=2242\", \"14\u00d725=350\"),\n    @(\"29\u00d712=348\", \"56\u00d779=4424\"),\n    @(\"82\u00d715=1230\", \"53\u00d713=689\"),\n    @(\"52\u00d744=2288\", \"14\u00d786=1204\"),\n    @(\"27\u00d729=783\", \"92\u00d775=6900\"),\n    @(\"24\u00d763=1512\", \"66\u00d723=1518\"),\n    @(\"86\u00d733=2838\", \"80\u00d742=3360\"),\n    @(\"83\u00d788=7304\", \"86\u00d759=5074\"),\n    @(\"56\u00d735=1960\", \"22\u00d776=1672\"),\n    @(\"28\u00d786=2408\", \"97\u00d744=4268\"),\n    @(\"87\u00d793=8091\", \"63\u00d754=3402\"),\n    @(\"30\u00d793=2790\", \"64\u00d733=2112\"),\n    @(\"46\u00d753=2438\", \"35\u00d744=1540\"),\n    @(\"33\u00d781=2673\", \"37\u00d725=925\"),\n    @(\"76\u00d799=7524\", \"17\u00d771=1207\"),\n    @(\"66\u00d779=5214\", \"76\u00d737=2812\"),\n    @(\"47\u00d716=752\", \"11\u00d756=616\"),\n    @(\"19\u00d756=1064\", \"63\u00d799=6237\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
